$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 33, shifting existing rows 33:43 down to 34:44
$ws.Rows.Item(33).Insert()

# Populate the newly inserted row 33 with the new data record
$ws.Range("A33").Value = 7
$ws.Range("B33").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C33").Value = "Ñuble"
$ws.Range("D33").Value = (Get-Date -Year 2022 -Month 8 -Day 9 -Hour 0 -Minute 0 -Second 0)
$ws.Range("E33").Value = 16
$ws.Range("F33").Value = 100112013
$ws.Range("G33").Value = "Alcachofa"
$ws.Range("H33").Value = "Madrigal"
$ws.Range("I33").Value = "Primera"
$ws.Range("J33").Value = 60
$ws.Range("K33").Value = 13000
$ws.Range("L33").Value = 14000
$ws.Range("M33").Value = 13500
$ws.Range("N33").Value = "$/caja 40 unidades"
$ws.Range("O33").Value = "Provincia de Limarí"
$ws.Range("P33").Value = 338
$ws.Range("Q33").Value = 40
$ws.Range("R33").Value = "Hortaliza"
